$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("M1").Value = "test"
